$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($r = 2; $r -le 153; $r++) {
    $ws.Range("C" + $r).Value = 45180
}
